$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that changes from
# 45203 (2023-10-04) to 45205 (2023-10-06) for every data row (2..351).
$startRow = 2
$endRow = 351

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45205
    }
}
